# Table 2 - completed summary statistics
# Adds two new worksheets ("Table 2 - city" and "Table 2 - country") after
# the existing "Table 2 - circle" sheet, each summarising the organic-vendor
# stats at a coarser grain (per city, per country).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 4: "Table 2 - city"
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$cityWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$cityWs.Name = "Table 2 - city"

# All values on this sheet are textual (e.g. "20 (33%)", "33%", "1.00"),
# matching the look of the sibling "Table 2 - circle" sheet. Force the
# whole used range to Text format *before* writing so Excel doesn't
# auto-convert look-alike numbers/percentages into real numbers.
$cityRange = $cityWs.Range("A1:G10")
$cityRange.NumberFormat = "@"

$cityData = @(
    @("country", "city", "vendors", "count_org_vendors", "org_vendor_perc", "median_org_foods_count", "iqr_org_foods_count"),
    @("Brazil", "Rio de Janeiro", "60", "20 (33%)", "33%", "0 (1.00)", "1.00"),
    @("Brazil", "Sao Paolo", "58", "4 (7%)", "7%", "0 (0.00)", "0.00"),
    @("Brazil", "Sinop", "8", "2 (25%)", "25%", "0 (0.50)", "0.50"),
    @("India", "Hyderabad", "141", "31 (22%)", "22%", "0 (0.00)", "0.00"),
    @("India", "Latur", "120", "71 (59%)", "59%", "1 (1.00)", "1.00"),
    @("India", "Visakhapatnam", "226", "53 (23%)", "23%", "0 (0.00)", "0.00"),
    @("UK", "Birmingham", "43", "17 (40%)", "40%", "0 (1.00)", "1.00"),
    @("UK", "Edinburgh", "60", "32 (53%)", "53%", "1 (2.25)", "2.25"),
    @("UK", "London", "92", "66 (72%)", "72%", "1 (2.00)", "2.00")
)

for ($r = 0; $r -lt $cityData.Length; $r++) {
    $row = $cityData[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $cityWs.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

$cityWs.PageSetup.PaperSize = 9
$cityWs.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# Sheet 5: "Table 2 - country"
# ---------------------------------------------------------------------
$countryWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $cityWs)
$countryWs.Name = "Table 2 - country"

$countryWs.Range("A1").Value = "country"
$countryWs.Range("B1").Value = "vendors"
$countryWs.Range("C1").Value = "count_org_vendors"
$countryWs.Range("D1").Value = "org_vendor_perc"
$countryWs.Range("E1").Value = "median_org_foods_count"
$countryWs.Range("F1").Value = "iqr_org_foods_count"

$countryWs.Range("A2").Value = "Brazil"
$countryWs.Range("B2").Value = 126
$countryWs.Range("C2").Value = 26
$countryWs.Range("D2").Value = 0.206349206349206
$countryWs.Range("E2").Value = 0
$countryWs.Range("F2").Value = 0

$countryWs.Range("A3").Value = "India"
$countryWs.Range("B3").Value = 487
$countryWs.Range("C3").Value = 155
$countryWs.Range("D3").Value = 0.318275154004107
$countryWs.Range("E3").Value = 0
$countryWs.Range("F3").Value = 1

$countryWs.Range("A4").Value = "UK"
$countryWs.Range("B4").Value = 195
$countryWs.Range("C4").Value = 115
$countryWs.Range("D4").Value = 0.58974358974359
$countryWs.Range("E4").Value = 1
$countryWs.Range("F4").Value = 2

$countryWs.PageSetup.PaperSize = 9
$countryWs.PageSetup.Orientation = 1

# Restore the original active sheet/selection so the new sheets don't end
# up "stealing" tab focus.
$wb.Worksheets.Item(1).Activate()
